# Figure3.pptx update: "update scripts for new data"
#
# - Merge "Tropical SW " + "Pacific" runs -> "Tropical SW Pacific"
# - Merge " " + "sites" runs -> " sites" (two panels: 9 sites / 6 sites)
# - Fix station count typo: "48 stations" -> "18 stations"
# - Drop a stray/redundant endParaRPr on the sub-figure label "b"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate shapes by name (robust to any re-ordering of Shapes.Item indices)
function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
    }
    return $null
}

# --- 1) "Tropical SW " + "Pacific" -> single run "Tropical SW Pacific" ---
$shpPacific = Get-ShapeByName $s "ZoneTexte 17"
$paraPacific = $shpPacific.TextFrame.TextRange.Paragraphs(1)
# Force a full-run rewrite so the two runs collapse into one.
$paraPacific.Text = "x"
$paraPacific.Text = "Tropical SW Pacific"

# --- 2) "9" + " " + "sites" -> "9" + " sites" ---
$shp9sites = Get-ShapeByName $s "ZoneTexte 22"
$para9 = $shp9sites.TextFrame.TextRange.Paragraphs(1)
$sub9 = $para9.Characters(2, $para9.Length - 1)
$sub9.Text = " sites"

# --- 3) "48 stations" -> "18 stations" (keep as 3 runs: "1" / "8 " / "stations") ---
$shpStations = Get-ShapeByName $s "ZoneTexte 23"
$paraStations = $shpStations.TextFrame.TextRange.Paragraphs(2)
$firstChar = $paraStations.Characters(1, 1)
$firstChar.Text = "1"
$tail = $paraStations.Characters(4, $paraStations.Length - 3)
$tail.Text = "x"
$tailRestore = $paraStations.Characters(4, 1)
$tailRestore.Text = "stations"

# --- 4) "6" + " " + "sites" -> "6" + " sites" ---
$shp6sites = Get-ShapeByName $s "ZoneTexte 25"
$para6 = $shp6sites.TextFrame.TextRange.Paragraphs(1)
$sub6 = $para6.Characters(2, $para6.Length - 1)
$sub6.Text = " sites"

# --- 5) Drop the redundant endParaRPr after the "b" sub-figure label ---
$shpB = Get-ShapeByName $s "ZoneTexte 57"
$trB = $shpB.TextFrame.TextRange
$trB.Text = "b" + [char]13 + "x"
$trB2 = $shpB.TextFrame.TextRange
$trB2.Text = "b"
